$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.017884162945029
$ws.Range("D2").Value = 4.139476750494767
$ws.Range("E2").Value = 16.51856462408147
$ws.Range("F2").Value = 24.3879733745967
$ws.Range("G2").Value = 29.77040432410453
$ws.Range("H2").Value = 13.90682887576727
$ws.Range("I2").Value = 20.07645073494477
$ws.Range("K2").Value = 14.35047164136688
$ws.Range("N2").Value = 16.33815385037597
$ws.Range("C3").Value = 2.954734837776295
$ws.Range("D3").Value = 4.146859750410423
$ws.Range("E3").Value = 15.57763463566987
$ws.Range("F3").Value = 24.22604852076245
$ws.Range("G3").Value = 29.38688244581259
$ws.Range("H3").Value = 13.92739512934794
$ws.Range("I3").Value = 20.01473977122937
$ws.Range("K3").Value = 13.74273257443791
$ws.Range("N3").Value = 16.39742679596634
$ws.Range("C4").Value = 2.914836956626297
$ws.Range("D4").Value = 4.15166696421736
$ws.Range("E4").Value = 14.97529892405943
$ws.Range("F4").Value = 24.13602020918362
$ws.Range("G4").Value = 29.16393252068305
$ws.Range("H4").Value = 13.94396108935942
$ws.Range("I4").Value = 19.98414529604162
$ws.Range("K4").Value = 13.35792127793047
$ws.Range("N4").Value = 16.43572990865561
$ws.Range("C5").Value = 2.89830868050164
$ws.Range("D5").Value = 4.153695174016311
$ws.Range("E5").Value = 14.72392582456065
$ws.Range("F5").Value = 24.10172461283616
$ws.Range("G5").Value = 29.07634777207079
$ws.Range("H5").Value = 13.95169747574747
$ws.Range("I5").Value = 19.97351609921173
$ws.Range("K5").Value = 13.19841104092379
$ws.Range("N5").Value = 16.45182016448745
$ws.Range("C6").Value = 2.895548287573455
$ws.Range("D6").Value = 4.154036147869497
$ws.Range("E6").Value = 14.68183710889577
$ws.Range("F6").Value = 24.09617506389384
$ws.Range("G6").Value = 29.06200507575881
$ws.Range("H6").Value = 13.95304148141506
$ws.Range("I6").Value = 19.97186220973818
$ws.Range("K6").Value = 13.17176906400334
$ws.Range("N6").Value = 16.45452105542148
$ws.Range("C7").Value = 2.914615124208556
$ws.Range("D7").Value = 4.151694036649781
$ws.Range("E7").Value = 14.97193237820357
$ws.Range("F7").Value = 24.13554796814707
$ws.Range("G7").Value = 29.16273793788506
$ws.Range("H7").Value = 13.94406144079614
$ws.Range("I7").Value = 19.98399450117478
$ws.Range("K7").Value = 13.35578065338854
$ws.Range("N7").Value = 16.43594495647426
$ws.Range("C8").Value = 2.996351154558705
$ws.Range("D8").Value = 4.141965728015183
$ws.Range("E8").Value = 16.19938394351127
$ws.Range("F8").Value = 24.33021114717113
$ws.Range("G8").Value = 29.63563789389922
$ws.Range("H8").Value = 13.91310006190693
$ws.Range("I8").Value = 20.05366141016674
$ws.Range("K8").Value = 14.1434745193367
$ws.Range("N8").Value = 16.35819568134732
$ws.Range("C9").Value = 3.147182123539721
$ws.Range("D9").Value = 4.125048797712806
$ws.Range("E9").Value = 18.49712254302333
$ws.Range("F9").Value = 24.78494356627702
$ws.Range("G9").Value = 30.65683586411075
$ws.Range("H9").Value = 13.88382354759696
$ws.Range("I9").Value = 20.24788801123382
$ws.Range("K9").Value = 15.58676535659887
$ws.Range("N9").Value = 16.22081931332749
$ws.Range("C10").Value = 3.251573021785089
$ws.Range("D10").Value = 4.113918578857382
$ws.Range("E10").Value = 20.14994154960957
$ws.Range("F10").Value = 25.16114363099285
$ws.Range("G10").Value = 31.45629704209972
$ws.Range("H10").Value = 13.88171388906335
$ws.Range("I10").Value = 20.4251490463236
$ws.Range("K10").Value = 16.57528513106083
$ws.Range("N10").Value = 16.12900560887659
$ws.Range("C11").Value = 3.297543512891347
$ws.Range("D11").Value = 4.109133527624878
$ws.Range("E11").Value = 20.8598812657774
$ws.Range("F11").Value = 25.34083206274394
$ws.Range("G11").Value = 31.82886908490165
$ws.Range("H11").Value = 13.88500684192414
$ws.Range("I11").Value = 20.51312086322277
$ws.Range("K11").Value = 17.00772167374474
$ws.Range("N11").Value = 16.08919973032446
$ws.Range("C12").Value = 3.314723716958916
$ws.Range("D12").Value = 4.107361284884194
$ws.Range("E12").Value = 21.12272506861252
$ws.Range("F12").Value = 25.41004984365883
$ws.Range("G12").Value = 31.97107741309905
$ws.Range("H12").Value = 13.8868677469344
$ws.Range("I12").Value = 20.54746971199597
$ws.Range("K12").Value = 17.16888220257952
$ws.Range("N12").Value = 16.0744069834789
$ws.Range("C13").Value = 3.311033930224363
$ws.Range("D13").Value = 4.107741204789379
$ws.Range("E13").Value = 21.06638301903524
$ws.Range("F13").Value = 25.39509128332156
$ws.Range("G13").Value = 31.94040281671155
$ws.Range("H13").Value = 13.88643963012518
$ws.Range("I13").Value = 20.54002633005666
$ws.Range("K13").Value = 17.13429038924773
$ws.Range("N13").Value = 16.07758039070477
$ws.Range("C14").Value = 3.298961556939775
$ws.Range("D14").Value = 4.108986928799855
$ws.Range("E14").Value = 20.88162552958837
$ws.Range("F14").Value = 25.34650346174143
$ws.Range("G14").Value = 31.84054702570746
$ws.Range("H14").Value = 13.88514762363196
$ws.Range("I14").Value = 20.51592612079096
$ws.Range("K14").Value = 17.02103302189377
$ws.Range("N14").Value = 16.08797709893876
$ws.Range("C15").Value = 3.291536937889226
$ws.Range("D15").Value = 4.109755140653832
$ws.Range("E15").Value = 20.76767655846357
$ws.Range("F15").Value = 25.31689314483618
$ws.Range("G15").Value = 31.77952415374072
$ws.Range("H15").Value = 13.88443624719684
$ws.Range("I15").Value = 20.50129834088937
$ws.Range("K15").Value = 16.95131874145608
$ws.Range("N15").Value = 16.09438193027399
$ws.Range("C16").Value = 3.248537367859772
$ws.Range("D16").Value = 4.114236868916734
$ws.Range("E16").Value = 20.10270448982859
$ws.Range("F16").Value = 25.14956794000153
$ws.Range("G16").Value = 31.432114253991
$ws.Range("H16").Value = 13.8815844791666
$ws.Range("I16").Value = 20.41954591501934
$ws.Range("K16").Value = 16.5466672178215
$ws.Range("N16").Value = 16.13164637813325
$ws.Range("C17").Value = 3.221762786736706
$ws.Range("D17").Value = 4.117057324292063
$ws.Range("E17").Value = 19.68405009861495
$ws.Range("F17").Value = 25.04907095728882
$ws.Range("G17").Value = 31.221156822574
$ws.Range("H17").Value = 13.88092610036542
$ws.Range("I17").Value = 20.3712593896936
$ws.Range("K17").Value = 16.29392282837038
$ws.Range("N17").Value = 16.15500830959167
$ws.Range("C18").Value = 3.206220452334477
$ws.Range("D18").Value = 4.118705769260973
$ws.Range("E18").Value = 19.43930185305858
$ws.Range("F18").Value = 24.99207638949454
$ws.Range("G18").Value = 31.100666127403
$ws.Range("H18").Value = 13.88094758709143
$ws.Range("I18").Value = 20.34417793971826
$ws.Range("K18").Value = 16.14693326871256
$ws.Range("N18").Value = 16.16863008040099
$ws.Range("C19").Value = 3.200933949570911
$ws.Range("D19").Value = 4.119268410846476
$ws.Range("E19").Value = 19.35575492354128
$ws.Range("F19").Value = 24.97291955674794
$ws.Range("G19").Value = 31.06002015553362
$ws.Range("H19").Value = 13.8810235092749
$ws.Range("I19").Value = 20.33512796037156
$ws.Range("K19").Value = 16.09689105163423
$ws.Range("N19").Value = 16.17327391702399
$ws.Range("C20").Value = 3.224627786537142
$ws.Range("D20").Value = 4.116754372941312
$ws.Range("E20").Value = 19.72902505422268
$ws.Range("F20").Value = 25.05968573996081
$ws.Range("G20").Value = 31.24352712372657
$ws.Range("H20").Value = 13.88095475099778
$ws.Range("I20").Value = 20.37632811601052
$ws.Range("K20").Value = 16.32099620315579
$ws.Range("N20").Value = 16.15250229401705
$ws.Range("C21").Value = 3.302513761568407
$ws.Range("D21").Value = 4.108619952463638
$ws.Range("E21").Value = 20.93605570335108
$ws.Range("F21").Value = 25.36074349243805
$ws.Range("G21").Value = 31.86984784251155
$ws.Range("H21").Value = 13.88551043868108
$ws.Range("I21").Value = 20.52297698668166
$ws.Range("K21").Value = 17.05437067290465
$ws.Range("N21").Value = 16.08491571911965
$ws.Range("C22").Value = 3.352084886597015
$ws.Range("D22").Value = 4.10353524263888
$ws.Range("E22").Value = 21.68999307431081
$ws.Range("F22").Value = 25.56431514713199
$ws.Range("G22").Value = 32.28566079229159
$ws.Range("H22").Value = 13.89206706533264
$ws.Range("I22").Value = 20.62484771778037
$ws.Range("K22").Value = 17.51851355757885
$ws.Range("N22").Value = 16.04238067411957
$ws.Range("C23").Value = 3.325752692966418
$ws.Range("D23").Value = 4.106227932934595
$ws.Range("E23").Value = 21.29078668915639
$ws.Range("F23").Value = 25.45506104856651
$ws.Range("G23").Value = 32.06319230524212
$ws.Range("H23").Value = 13.88823951417854
$ws.Range("I23").Value = 20.56993274909302
$ws.Range("K23").Value = 17.27221212182573
$ws.Range("N23").Value = 16.0649330128786
$ws.Range("C24").Value = 3.223332984327677
$ws.Range("D24").Value = 4.116891253319636
$ws.Range("E24").Value = 19.70870452053538
$ws.Range("F24").Value = 25.05488435505783
$ws.Range("G24").Value = 31.23341103313186
$ws.Range("H24").Value = 13.88094055223631
$ws.Range("I24").Value = 20.37403442801015
$ws.Range("K24").Value = 16.30876156513839
$ws.Range("N24").Value = 16.15363466958244
$ws.Range("C25").Value = 3.10746169333677
$ws.Range("D25").Value = 4.129395950670899
$ws.Range("E25").Value = 17.85103401436508
$ws.Range("F25").Value = 24.65434109178642
$ws.Range("G25").Value = 30.37133471975421
$ws.Range("H25").Value = 13.88835186453732
$ws.Range("I25").Value = 20.18922408677951
$ws.Range("K25").Value = 15.20830292669651
$ws.Range("N25").Value = 16.25637644268217
